# Konklusion.docx revision pass (Anders Meidahl, 2016-05-18)
$d = $word.ActiveDocument

$word.UserName = "Anders Meidahl"
$word.UserInitials = "AM"

# ---------------------------------------------------------------------
# 1) Comment on the "Konklusion" heading.
# ---------------------------------------------------------------------
$headingRange = $d.Paragraphs(1).Range
$headingRange = $d.Range($headingRange.Start, $headingRange.End - 1)
$d.Comments.Add($headingRange, "Mere kød og evt noget med vores brugerundersøgelse") | Out-Null

# ---------------------------------------------------------------------
# 2) "...at handle i så billigt..." -> delete "i " (tracked deletion).
# ---------------------------------------------------------------------
$d.TrackRevisions = $true

$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Execute("handle i så billigt", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$delRange = $d.Range($rng.Start + 7, $rng.Start + 9)
$delRange.Delete()

# ---------------------------------------------------------------------
# 3) Paragraph 3: drop three spaces ("computer applikationer" ->
#    "computerapplikationer", "administrativ applikation" ->
#    "administrativapplikation", "forbruger applikation" ->
#    "forbrugerapplikation"), and move the _GoBack bookmark.
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("computer applikationer", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$spacePos = $rng.Start + [string]"computer".Length
$d.Range($spacePos, $spacePos + 1).Delete()

$rng = $d.Content
$rng.Find.Execute("administrativ applikation", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$spacePos = $rng.Start + [string]"administrativ".Length
$d.Range($spacePos, $spacePos + 1).Delete()

$rng = $d.Content
$rng.Find.Execute("forbruger applikation", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$spacePos = $rng.Start + [string]"forbruger".Length
$d.Range($spacePos, $spacePos + 1).Delete()

# Move the (hidden) _GoBack bookmark to just after "...sin indkøbsliste til".
$rng = $d.Content
$rng.Find.Execute("sende sin indkøbsliste til sin mail", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bmPos = $rng.Start + [string]"sende sin indkøbsliste til".Length
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# ---------------------------------------------------------------------
# 4) "...I databasen...Forbruger applikationen kan ikke tilføje noget
#    til information..." -> drop the space after "Forbruger" and drop
#    "noget til ".
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Forbruger applikationen kan ikke tilføje noget til information", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$spacePos = $rng.Start + [string]"Forbruger".Length
$d.Range($spacePos, $spacePos + 1).Delete()

$rng = $d.Content
$rng.Find.Execute("tilføje noget til information", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$dropStart = $rng.Start + [string]"tilføje ".Length
$dropEnd = $dropStart + [string]"noget til ".Length
$d.Range($dropStart, $dropEnd).Delete()

# ---------------------------------------------------------------------
# 5) "...potentialet for at blive..." -> "for " deleted, "til " inserted.
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("potentialet for at blive", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$forStart = $rng.Start + [string]"potentialet ".Length
$forEnd = $forStart + [string]"for ".Length
$forRange = $d.Range($forStart, $forEnd)
$forRange.Delete()
$insPoint = $d.Range($forStart, $forStart)
$insPoint.InsertBefore("til ")

$d.TrackRevisions = $false
